# "Generate Report for Handback" — fills in the handback columns for the
# 8701aa89-d787-4808-a983-8437cf16ed01 row (row 5) on the zh-cn and de-de
# localization-status sheets, now that the handback has been processed.

$wb = $excel.ActiveWorkbook

$fileDisplay = "8701aa89-d787-4808-a983-8437cf16ed01.md"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b950b4fadf73b346ffc4a8cb05eaf97231aa0570/e2e/8701aa89-d787-4808-a983-8437cf16ed01.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c90368aadecd20bad147eb4adbfe359feb2a81c9/e2e/8701aa89-d787-4808-a983-8437cf16ed01.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(16).ColumnWidth = 40

$wsZh.Range("J5").Value = "8701aa89-d787-4808-a983-8437cf16ed01.a0ea39dbd133e8bc91eb581552a1cc6a71d95845.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-09-07 05:58:15"
$wsZh.Range("P5").Value = $errorMsg

$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b950b4fadf73b346ffc4a8cb05eaf97231aa0570/e2e/8701aa89-d787-4808-a983-8437cf16ed01.md", [Type]::Missing, [Type]::Missing, $fileDisplay) | Out-Null

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(16).ColumnWidth = 40

$wsDe.Range("J5").Value = "8701aa89-d787-4808-a983-8437cf16ed01.a0ea39dbd133e8bc91eb581552a1cc6a71d95845.de-de.xlf"
$wsDe.Range("K5").Value = "2016-09-07 05:58:33"
$wsDe.Range("P5").Value = $errorMsg

$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b950b4fadf73b346ffc4a8cb05eaf97231aa0570/e2e/8701aa89-d787-4808-a983-8437cf16ed01.md", [Type]::Missing, [Type]::Missing, $fileDisplay) | Out-Null
